# Auto-generated edit script: update cryptos list Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.048.09"
$ws.Range("D3").Value = "3.590.12"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'581.05"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'191.19"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("D8").Value = "3.586.52"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'0.182"
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "'55.87"
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("E13").Value = "  +6.95%  "
$ws.Range("D14").Value = "'9.68"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "4.175.16"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'20.03"
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("D17").Value = "3.595.22"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "70.084.58"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "'480.31"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "'19.17"
$ws.Range("E23").Value = "  +10.35%  "
$ws.Range("D24").Value = "'5.03"
$ws.Range("E24").Value = "  -6.32%  "
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "'95.57"
$ws.Range("E26").Value = "  +6.04%  "
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "'11.12"
$ws.Range("D29").Value = "'9.42"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "'32.17"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'7.66"
$ws.Range("E31").Value = "  +3.03%  "
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").Value = "'12.23"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "'66.68"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").Value = "'583.88"
$ws.Range("E35").Value = "  -6.11%  "
$ws.Range("D36").Value = "'39.03"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "0.0₃0804"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").Value = "'0.397"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "'3.30"
$ws.Range("E40").Value = "  +23.71%  "
$ws.Range("D41").Value = "'3.46"
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("D42").Value = "3.222.81"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("E44").Value = "  +7.10%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").Value = "'9.50"
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "'3.14"
$ws.Range("E51").Value = "  -5.39%  "
